$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mismatched profile values
$ws.Range("A2").Value = "JImport1CK"
$ws.Range("A3").Value = "Jimport2CK2017"

# Add new "Orin" column with values
$ws.Range("H1").Value = "Orin"
$ws.Range("H2").Value = 20189
$ws.Range("H3").Value = 20199

$ws.Range("A2").Select()
